$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Roraima
$ws.Range("A2").Value = "Roraima"
$ws.Range("B2").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C2").Value = 2.4
$ws.Range("D2").Value = "1º"

# Row 3: Piauí (was Amapá)
$ws.Range("A3").Value = "Piauí"
$ws.Range("B3").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C3").Value = 1.1
$ws.Range("D3").Value = "2º"

# Row 4: Ceará (was Piauí)
$ws.Range("A4").Value = "Ceará"
$ws.Range("B4").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C4").Value = 0.8999999999999995
$ws.Range("D4").Value = "3º"

# Row 5: Amapá (was Ceará)
$ws.Range("A5").Value = "Amapá"
$ws.Range("B5").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C5").Value = 0.8999999999999986
$ws.Range("D5").Value = "4º"

# Row 6: Mato Grosso do Sul (was Amazonas)
$ws.Range("A6").Value = "Mato Grosso do Sul"
$ws.Range("B6").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C6").Value = 0.7000000000000002
$ws.Range("D6").Value = "5º"

# Row 7: Rondônia (was Goiás)
$ws.Range("A7").Value = "Rondônia"
$ws.Range("B7").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C7").Value = 0.6999999999999997
$ws.Range("D7").Value = "6º"

# Row 8: Sergipe
$ws.Range("A8").Value = "Sergipe"
$ws.Range("B8").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C8").Value = -0.7000000000000011
$ws.Range("D8").Value = "20º"

# Row 9: Nordeste
$ws.Range("A9").Value = "Nordeste"
$ws.Range("B9").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C9").Value = -0.5

# Row 10: Brasil
$ws.Range("A10").Value = "Brasil"
$ws.Range("B10").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C10").Value = -0.5
